$d = $word.ActiveDocument

# 1. Heading "Datos del estudiante" -> "Datos del Alumno/a:"
$d.Content.Find.Execute("Datos del estudiante", $false, $false, $false, $false, $false, $true, 1, $false, "Datos del Alumno/a:", 2)

# 2. "Adultos autorizados a retirar al estudiante:" -> "Adultos autorizados a retirar al Alumno/a:"
$d.Content.Find.Execute("Adultos autorizados a retirar al estudiante:", $false, $false, $false, $false, $false, $true, 1, $false, "Adultos autorizados a retirar al Alumno/a:", 2)

# 3. "Ante una emergencia se trasladará al estudiante al hospital más cercano: " -> "... al alumna/o al hospital más cercano: "
$d.Content.Find.Execute("Ante una emergencia se trasladará al estudiante al hospital más cercano: ", $false, $false, $false, $false, $false, $true, 1, $false, "Ante una emergencia se trasladará al alumna/o al hospital más cercano: ", 2)

# 4. "El/la que suscribe, responsable del/ de la estudiante ..." -> "...responsable del/ de la alumna/o ..."
$d.Content.Find.Execute("El/la que suscribe, responsable del/ de la estudiante __________________________________ con DNI ____________________ del", $false, $false, $false, $false, $false, $true, 1, $false, "El/la que suscribe, responsable del/ de la alumna/o __________________________________ con DNI ____________________ del", 2)

# 5. " _________________ toma conocimiento y autoriza a que el/la estudiante realice actividad." -> "...el/la alumna/o realice actividad."
$d.Content.Find.Execute(" _________________ toma conocimiento y autoriza a que el/la estudiante realice actividad.", $false, $false, $false, $false, $false, $true, 1, $false, " _________________ toma conocimiento y autoriza a que el/la alumna/o realice actividad.", 2)
